# Update the "Generate Report for Handback" timestamps.
# Overview sheet: Latest HO Xliff Generate Date for the 25d627c2... row (row 3)
# zh-cn sheet: Correspond Handoff / Handback DateTime for the 25d627c2... row (row 3)
# de-de sheet: Correspond Handoff DateTime for the 25d627c2... row (row 3)
#   (the de-de Correspond Handoff DateTime cell previously shared the exact same
#   text as the Overview cell, so it must be updated in lockstep.)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Sheets.Item("Overview")
$wsZhCn     = $wb.Sheets.Item("zh-cn")
$wsDeDe     = $wb.Sheets.Item("de-de")

$wsOverview.Range("G3").Value = "2016-08-18 08:48:22"

$wsZhCn.Range("H3").Value = "2016-08-18 08:48:16"
$wsZhCn.Range("K3").Value = "2016-08-18 08:48:34"

$wsDeDe.Range("H3").Value = "2016-08-18 08:48:22"
$wsDeDe.Range("K3").Value = "2016-08-18 08:48:42"
